# Regenerate merged AHB files
# - Rename the "_old"/"_new" header-column suffixes to "_FV2310"/"_FV2404"
# - Freeze the header row
# - Turn the data range into a native Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (A1:J1 = *_old -> *_FV2310, L1:U1 = *_new -> *_FV2404) ---
$ws.Cells.Item(1, 1).Value  = "Segmentname_FV2310"
$ws.Cells.Item(1, 2).Value  = "Segmentgruppe_FV2310"
$ws.Cells.Item(1, 3).Value  = "Segment_FV2310"
$ws.Cells.Item(1, 4).Value  = "Datenelement_FV2310"
$ws.Cells.Item(1, 5).Value  = "Segment ID_FV2310"
$ws.Cells.Item(1, 6).Value  = "Code_FV2310"
$ws.Cells.Item(1, 7).Value  = "Qualifier_FV2310"
$ws.Cells.Item(1, 8).Value  = "Beschreibung_FV2310"
$ws.Cells.Item(1, 9).Value  = "Bedingungsausdruck_FV2310"
$ws.Cells.Item(1, 10).Value = "Bedingung_FV2310"
# column 11 ("diff") is unchanged

$ws.Cells.Item(1, 12).Value = "Segmentname_FV2404"
$ws.Cells.Item(1, 13).Value = "Segmentgruppe_FV2404"
$ws.Cells.Item(1, 14).Value = "Segment_FV2404"
$ws.Cells.Item(1, 15).Value = "Datenelement_FV2404"
$ws.Cells.Item(1, 16).Value = "Segment ID_FV2404"
$ws.Cells.Item(1, 17).Value = "Code_FV2404"
$ws.Cells.Item(1, 18).Value = "Qualifier_FV2404"
$ws.Cells.Item(1, 19).Value = "Beschreibung_FV2404"
$ws.Cells.Item(1, 20).Value = "Bedingungsausdruck_FV2404"
$ws.Cells.Item(1, 21).Value = "Bedingung_FV2404"

# --- 2. Freeze the top (header) row ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null

# --- 3. Convert the data range into a native table ---
$range = $ws.Range("A1:U89")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
